# Rebuild the slide with native, editable PPTX shapes (text boxes + rounded
# rectangles) in place of the old placeholder-text shapes + flattened
# background picture. IDs need to restart at 2 (matching the target OOXML),
# so we build everything on a brand-new slide (same layout as the original)
# and then drop the old slide.

$p = $ppt.ActivePresentation
$oldSlide = $p.Slides.Item(1)
$layout = $oldSlide.CustomLayout
$newSlide = $p.Slides.AddSlide(2, $layout)

# id=2 TextBox 1
$sh2 = $newSlide.Shapes.AddTextbox(1, 1321.0127559055118, 662.2127559055118, 262.0007874015748, 149.9975590551181)
$sh2.Name = "TextBox 1"
$sh2.Fill.Visible = $false
$sh2.TextFrame.WordWrap = $true
$sh2.TextFrame.AutoSize = 1
$sh2.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh2.TextFrame.TextRange.Font.Size = 26
$sh2.TextFrame.TextRange.Font.Bold = $false
$sh2.TextFrame.TextRange.Font.Italic = $false
$sh2.TextFrame.TextRange.Font.Color.RGB = 16777215
$sh2.TextFrame.TextRange.Text = "Composite trading result. `rShould we be invested?"
$sh2.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh2.TextFrame.TextRange.Paragraphs(2).Font.Name = "Quicksand (TT)"
$sh2.TextFrame.TextRange.Paragraphs(2).Font.Size = 26
$sh2.TextFrame.TextRange.Paragraphs(2).Font.Color.RGB = 16777215
$sh2.Left = 1321.0127559055118
$sh2.Top = 662.2127559055118
$sh2.Width = 262.0007874015748
$sh2.Height = 149.9975590551181

# id=3 TextBox 2
$sh3 = $newSlide.Shapes.AddTextbox(1, 1346.8751968503936, 500.3927559055118, 210.26874015748032, 30.00236220472441)
$sh3.Name = "TextBox 2"
$sh3.Fill.Visible = $false
$sh3.TextFrame.WordWrap = $true
$sh3.TextFrame.AutoSize = 1
$sh3.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh3.TextFrame.TextRange.Font.Size = 30
$sh3.TextFrame.TextRange.Font.Bold = $false
$sh3.TextFrame.TextRange.Font.Italic = $false
$sh3.TextFrame.TextRange.Font.Color.RGB = 2171169
$sh3.TextFrame.TextRange.Text = "Composites"
$sh3.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh3.Left = 1346.8751968503936
$sh3.Top = 500.3927559055118
$sh3.Width = 210.26874015748032
$sh3.Height = 30.00236220472441

# id=4 Rounded Rectangle 3
$sh4 = $newSlide.Shapes.AddShape(5, 1301.1623622047243, 481.9031496062992, 301.6871653543307, 72.6407874015748)
$sh4.Name = "Rounded Rectangle 3"
$sh4.Fill.ForeColor.RGB = 1032446
$sh4.Line.Visible = $false
$sh4.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh4.Left = 1301.1623622047243
$sh4.Top = 481.9031496062992
$sh4.Width = 301.6871653543307
$sh4.Height = 72.6407874015748

# id=5 Rounded Rectangle 4
$sh5 = $newSlide.Shapes.AddShape(5, 1299.146377952756, 597.6647244094488, 305.7263779527559, 279.0863779527559)
$sh5.Name = "Rounded Rectangle 4"
$sh5.Fill.ForeColor.RGB = 1032446
$sh5.Line.Visible = $false
$sh5.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh5.Left = 1299.146377952756
$sh5.Top = 597.6647244094488
$sh5.Width = 305.7263779527559
$sh5.Height = 279.0863779527559

# id=6 TextBox 5
$sh6 = $newSlide.Shapes.AddTextbox(1, 828.7559842519685, 662.2127559055118, 262.0007874015748, 149.9975590551181)
$sh6.Name = "TextBox 5"
$sh6.Fill.Visible = $false
$sh6.TextFrame.WordWrap = $true
$sh6.TextFrame.AutoSize = 1
$sh6.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh6.TextFrame.TextRange.Font.Size = 26
$sh6.TextFrame.TextRange.Font.Bold = $false
$sh6.TextFrame.TextRange.Font.Italic = $false
$sh6.TextFrame.TextRange.Font.Color.RGB = 16777215
$sh6.TextFrame.TextRange.Text = "Capture short-term market inefficiencies that generate high returns while invested."
$sh6.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh6.Left = 828.7559842519685
$sh6.Top = 662.2127559055118
$sh6.Width = 262.0007874015748
$sh6.Height = 149.9975590551181

# id=7 TextBox 6
$sh7 = $newSlide.Shapes.AddTextbox(1, 843.8975590551181, 500.3927559055118, 231.7175590551181, 30.00236220472441)
$sh7.Name = "TextBox 6"
$sh7.Fill.Visible = $false
$sh7.TextFrame.WordWrap = $true
$sh7.TextFrame.AutoSize = 1
$sh7.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh7.TextFrame.TextRange.Font.Size = 30
$sh7.TextFrame.TextRange.Font.Bold = $false
$sh7.TextFrame.TextRange.Font.Italic = $false
$sh7.TextFrame.TextRange.Font.Color.RGB = 2171169
$sh7.TextFrame.TextRange.Text = "Trigger Systems"
$sh7.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh7.Left = 843.8975590551181
$sh7.Top = 500.3927559055118
$sh7.Width = 231.7175590551181
$sh7.Height = 30.00236220472441

# id=8 Rounded Rectangle 7
$sh8 = $newSlide.Shapes.AddShape(5, 808.9127559055119, 481.9031496062992, 301.6871653543307, 72.6407874015748)
$sh8.Name = "Rounded Rectangle 7"
$sh8.Fill.ForeColor.RGB = 1032446
$sh8.Line.Visible = $false
$sh8.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh8.Left = 808.9127559055119
$sh8.Top = 481.9031496062992
$sh8.Width = 301.6871653543307
$sh8.Height = 72.6407874015748

# id=9 Rounded Rectangle 8
$sh9 = $newSlide.Shapes.AddShape(5, 806.8895275590551, 597.6647244094488, 305.7263779527559, 279.0863779527559)
$sh9.Name = "Rounded Rectangle 8"
$sh9.Fill.ForeColor.RGB = 1032446
$sh9.Line.Visible = $false
$sh9.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh9.Left = 806.8895275590551
$sh9.Top = 597.6647244094488
$sh9.Width = 305.7263779527559
$sh9.Height = 279.0863779527559

# id=10 TextBox 9
$sh10 = $newSlide.Shapes.AddTextbox(1, 336.9887401574803, 662.2127559055118, 262.0007874015748, 149.9975590551181)
$sh10.Name = "TextBox 9"
$sh10.Fill.Visible = $false
$sh10.TextFrame.WordWrap = $true
$sh10.TextFrame.AutoSize = 1
$sh10.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh10.TextFrame.TextRange.Font.Size = 26
$sh10.TextFrame.TextRange.Font.Bold = $false
$sh10.TextFrame.TextRange.Font.Italic = $false
$sh10.TextFrame.TextRange.Font.Color.RGB = 16777215
$sh10.TextFrame.TextRange.Text = "Total market systems that trade infrequently to capture long term trend changes."
$sh10.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh10.Left = 336.9887401574803
$sh10.Top = 662.2127559055118
$sh10.Width = 262.0007874015748
$sh10.Height = 149.9975590551181

# id=11 TextBox 10
$sh11 = $newSlide.Shapes.AddTextbox(1, 362.8583464566929, 500.3927559055118, 210.26874015748032, 30.00236220472441)
$sh11.Name = "TextBox 10"
$sh11.Fill.Visible = $false
$sh11.TextFrame.WordWrap = $true
$sh11.TextFrame.AutoSize = 1
$sh11.TextFrame.TextRange.Font.Name = "Quicksand (TT)"
$sh11.TextFrame.TextRange.Font.Size = 30
$sh11.TextFrame.TextRange.Font.Bold = $false
$sh11.TextFrame.TextRange.Font.Italic = $false
$sh11.TextFrame.TextRange.Font.Color.RGB = 2171169
$sh11.TextFrame.TextRange.Text = "Base Systems"
$sh11.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh11.Left = 362.8583464566929
$sh11.Top = 500.3927559055118
$sh11.Width = 210.26874015748032
$sh11.Height = 30.00236220472441

# id=12 Rounded Rectangle 11
$sh12 = $newSlide.Shapes.AddShape(5, 317.1455905511811, 481.9031496062992, 301.6871653543307, 72.6407874015748)
$sh12.Name = "Rounded Rectangle 11"
$sh12.Fill.ForeColor.RGB = 1032446
$sh12.Line.Visible = $false
$sh12.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh12.Left = 317.1455905511811
$sh12.Top = 481.9031496062992
$sh12.Width = 301.6871653543307
$sh12.Height = 72.6407874015748

# id=13 Rounded Rectangle 12
$sh13 = $newSlide.Shapes.AddShape(5, 315.12952755905513, 597.6647244094488, 305.7263779527559, 279.0863779527559)
$sh13.Name = "Rounded Rectangle 12"
$sh13.Fill.ForeColor.RGB = 1032446
$sh13.Line.Visible = $false
$sh13.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh13.Left = 315.12952755905513
$sh13.Top = 597.6647244094488
$sh13.Width = 305.7263779527559
$sh13.Height = 279.0863779527559

# id=14 TextBox 13
$sh14 = $newSlide.Shapes.AddTextbox(1, 1186.9127559055119, 704.663937007874, 37.93677165354331, 65.10236220472441)
$sh14.Name = "TextBox 13"
$sh14.Fill.Visible = $false
$sh14.TextFrame.WordWrap = $true
$sh14.TextFrame.AutoSize = 1
$sh14.TextFrame.TextRange.Font.Name = "Rajdhani"
$sh14.TextFrame.TextRange.Font.Size = 70
$sh14.TextFrame.TextRange.Font.Bold = $false
$sh14.TextFrame.TextRange.Font.Italic = $false
$sh14.TextFrame.TextRange.Font.Color.RGB = 707583
$sh14.TextFrame.TextRange.Text = "="
$sh14.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh14.Left = 1186.9127559055119
$sh14.Top = 704.663937007874
$sh14.Width = 37.93677165354331
$sh14.Height = 65.10236220472441

# id=15 TextBox 14
$sh15 = $newSlide.Shapes.AddTextbox(1, 695.1455905511812, 704.663937007874, 37.44716535433071, 65.10236220472441)
$sh15.Name = "TextBox 14"
$sh15.Fill.Visible = $false
$sh15.TextFrame.WordWrap = $true
$sh15.TextFrame.AutoSize = 1
$sh15.TextFrame.TextRange.Font.Name = "Rajdhani"
$sh15.TextFrame.TextRange.Font.Size = 70
$sh15.TextFrame.TextRange.Font.Bold = $false
$sh15.TextFrame.TextRange.Font.Italic = $false
$sh15.TextFrame.TextRange.Font.Color.RGB = 707583
$sh15.TextFrame.TextRange.Text = "+"
$sh15.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh15.Left = 695.1455905511812
$sh15.Top = 704.663937007874
$sh15.Width = 37.44716535433071
$sh15.Height = 65.10236220472441

# id=16 TextBox 15
$sh16 = $newSlide.Shapes.AddTextbox(1, 84.0023622047244, 90.0, 419.99755905511813, 44.02795275590551)
$sh16.Name = "TextBox 15"
$sh16.Fill.Visible = $false
$sh16.TextFrame.WordWrap = $true
$sh16.TextFrame.AutoSize = 1
$sh16.TextFrame.TextRange.Font.Name = "Rajdhani"
$sh16.TextFrame.TextRange.Font.Size = 36
$sh16.TextFrame.TextRange.Font.Bold = $false
$sh16.TextFrame.TextRange.Font.Italic = $false
$sh16.TextFrame.TextRange.Font.Color.RGB = 16777215
$sh16.TextFrame.TextRange.Text = "process"
$sh16.TextFrame.TextRange.ParagraphFormat.Alignment = 1
$sh16.Left = 84.0023622047244
$sh16.Top = 90.0
$sh16.Width = 419.99755905511813
$sh16.Height = 44.02795275590551

# id=17 TextBox 16
$sh17 = $newSlide.Shapes.AddTextbox(1, 303.45834645669294, 194.91834645669292, 1313.0855905511812, 127.8)
$sh17.Name = "TextBox 16"
$sh17.Fill.Visible = $false
$sh17.TextFrame.WordWrap = $true
$sh17.TextFrame.AutoSize = 1
$sh17.TextFrame.TextRange.Font.Name = "Rajdhani"
$sh17.TextFrame.TextRange.Font.Size = 60
$sh17.TextFrame.TextRange.Font.Bold = $false
$sh17.TextFrame.TextRange.Font.Italic = $false
$sh17.TextFrame.TextRange.Font.Color.RGB = 16777215
$sh17.TextFrame.TextRange.Text = "Our Composites are Designed to Highlight times of Risk-On and Risk-Off Behavior"
$sh17.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$sh17.Left = 303.45834645669294
$sh17.Top = 194.91834645669292
$sh17.Width = 1313.0855905511812
$sh17.Height = 127.8

$oldSlide.Delete()

